# Applies the README/report edits described in the commit:
#   - rewrite the intro sentence (para 1)
#   - rewrite 'For this assignment ... features' (para 4)
#   - drop stray grammar-check markers with no text change (paras 5, 19)
#
# Each paragraph's full markup is replaced in one shot via Range.InsertXML
# (WordprocessingML package fragment) so the run layout matches the target
# exactly, instead of relying on Find&Replace's automatic run-merging.
$d = $word.ActiveDocument

# --- Paragraph 1 -- intro sentence rewrite ---
$xmlpara1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="43072C83" w14:textId="3231212B" w:rsidR="0014507E" w:rsidRDefault="007E3BD1"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="006F7DD0"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>In this video I will give you a brief description about what are chatbot is based on</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>what it does</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> and the features </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>based on</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> the </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>API’s we</w:t></w:r><w:r w:rsidRPr="006F7DD0"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> have implemented that would enhance the functionality of our bot.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void] $d.Paragraphs(1).Range.InsertXML($xmlpara1)

# --- Paragraph 4 -- 'For this assignment ... features' rewrite ---
$xmlpara4 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="6BE365E8" w14:textId="77777777" w:rsidR="004B6B91" w:rsidRDefault="006F7DD0"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">For this assignment </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>I have</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> decided to implement</w:t></w:r><w:r w:rsidR="00A41B5E"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">more </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>features</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void] $d.Paragraphs(4).Range.InsertXML($xmlpara4)

# --- Paragraph 5 -- drop proofErr markers around 'doesn't' / 'in itself is' (no text change) ---
$xmlpara5 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="2600CD4D" w14:textId="5AE89212" w:rsidR="00921DC9" w:rsidRDefault="00921DC9" w:rsidP="004B6B91"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Pos tagging: The Pos tagging feature allows the chatbot to recognize and label different parts of speech as nouns, pronouns adjectives etc.</w:t></w:r><w:r w:rsidR="00AE098C"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> We made use of POS tagging to recognize proper nouns, so that if there is a question in a particular topic that the user asks the bot about, and the bot doesn’t know the answer to it than the bot can retrieve the information from Wikipedia and give it to the user.</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> We have</w:t></w:r><w:r w:rsidR="00AE098C"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> also</w:t></w:r><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> used POS tagging in conjunction with sentiment analysis to judge the different emotions of the user more </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>accuratelty</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>. Since sentiment analysis in itself is not always accurate, so we had to use POS to make our bot more precise with its judgements.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void] $d.Paragraphs(5).Range.InsertXML($xmlpara5)

# --- Paragraph 19 -- drop proofErr markers around 'and also' (no text change) ---
$xmlpara19 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="3625CCCF" w14:textId="25FEA262" w:rsidR="003B065E" w:rsidRDefault="003B065E" w:rsidP="00DE36FE"><w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Hung: he was responsible for </w:t></w:r><w:r w:rsidR="000F3A42"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">updating and the WBS and the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="000F3A42"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>gantt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="000F3A42"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> chart and also fixed the bugs</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
[void] $d.Paragraphs(19).Range.InsertXML($xmlpara19)

Write-Output "edit complete"
